$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct the two mistranslated / misspelled strings in the language table
$ws.Range("B21").Value = "Could not connect to server"
$ws.Range("B10").Value = "Could not share"

# Update the view selection to match the author's final cursor position
$ws.Range("B10").Select()
